# Upload new version with timestamp
# Refresh the item-shortage report row with the latest snapshot figures:
#  - H7 (الرصيد الحالي / current balance) updates from "0:1" to "0:0"
#  - Q7 (عدد التعاملات / transaction count) updates from "0:0" to "0:1"
#  - P7 (سعر البيع / selling price) updates from "0.0000" to "20.0000"
#  - N8 (summary/total cell) updates from 0 to 20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# P7 must stay a literal text value (it already was text, formatted as "0.00"),
# so temporarily force a text number-format while assigning the new value,
# then restore the original number format to avoid altering the cell style.
$p7Range = $ws.Range("P7")
$p7Format = $p7Range.NumberFormat
$p7Range.NumberFormat = "@"
$p7Range.Value = "20.0000"
$p7Range.NumberFormat = $p7Format

# H7 and Q7 simply swap their text contents.
$ws.Range("H7").Value = "0:0"
$ws.Range("Q7").Value = "0:1"

# N8 total cell becomes numeric 20.
$ws.Range("N8").Value = 20
